$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking values are not
# auto-converted to floating point numbers by the Value setter.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.785.85"
$ws.Range("E2").Value = "  -4.69%  "

$ws.Range("D3").Value = "2.991.43"
$ws.Range("E3").Value = "  -4.72%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "542.65"
$ws.Range("E5").Value = "  -5.52%  "

$ws.Range("D6").Value = "151.94"
$ws.Range("E6").Value = "  -7.76%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  -0.75%  "

$ws.Range("D9").Value = "3.006.10"
$ws.Range("E9").Value = "  -4.69%  "

$ws.Range("E10").Value = "  -3.99%  "

$ws.Range("D11").Value = "6.15"
$ws.Range("E11").Value = "  -7.36%  "

$ws.Range("E12").Value = "  -3.36%  "

$ws.Range("D13").Value = "3.511.53"
$ws.Range("E13").Value = "  -4.74%  "

$ws.Range("E14").Value = "  -1.18%  "

$ws.Range("D15").Value = "61.865.47"
$ws.Range("E15").Value = "  -4.59%  "

$ws.Range("D16").Value = "24.04"
$ws.Range("E16").Value = "  -3.96%  "

$ws.Range("D17").Value = "2.998.00"
$ws.Range("E17").Value = "  -4.74%  "

$ws.Range("E18").Value = "  -5.73%  "

$ws.Range("D19").Value = "5.17"
$ws.Range("E19").Value = "  -1.36%  "

$ws.Range("D20").Value = "12.10"
$ws.Range("E20").Value = "  -3.20%  "

$ws.Range("D21").Value = "379.18"
$ws.Range("E21").Value = "  -8.46%  "

$ws.Range("D22").Value = "6.73"
$ws.Range("E22").Value = "  -4.43%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("E24").Value = "  -3.83%  "

$ws.Range("D25").Value = "66.14"
$ws.Range("E25").Value = "  -3.92%  "

$ws.Range("D26").Value = "3.114.44"
$ws.Range("E26").Value = "  -4.78%  "

$ws.Range("E27").Value = "  -2.72%  "

$ws.Range("E28").Value = "  -2.82%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").Value = "0.0₃0935"
$ws.Range("E30").Value = "  -9.66%  "

$ws.Range("D31").Value = "8.25"
$ws.Range("E31").Value = "  -9.38%  "

$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").Value = "1.73"
$ws.Range("E33").Value = "  -4.60%  "

$ws.Range("D34").Value = "20.49"
$ws.Range("E34").Value = "  -3.86%  "

$ws.Range("D35").Value = "161.14"
$ws.Range("E35").Value = "  -1.30%  "

$ws.Range("D36").Value = "6.00"
$ws.Range("E36").Value = "  -4.11%  "

$ws.Range("E37").Value = "  -5.77%  "

$ws.Range("E38").Value = "  -5.07%  "

$ws.Range("E39").Value = "  -5.70%  "

$ws.Range("D40").Value = "1.56"
$ws.Range("E40").Value = "  -7.74%  "

$ws.Range("E41").Value = "  -1.80%  "

$ws.Range("D42").Value = "2.419.33"
$ws.Range("E42").Value = "  -7.79%  "

$ws.Range("D43").Value = "3.91"
$ws.Range("E43").Value = "  -5.92%  "

$ws.Range("D44").Value = "22.13"
$ws.Range("E44").Value = "  -7.11%  "

$ws.Range("D45").Value = "0.673"
$ws.Range("E45").Value = "  -2.82%  "

$ws.Range("D46").Value = "0.0593"
$ws.Range("E46").Value = "  -3.85%  "

$ws.Range("E47").Value = "  -2.69%  "

$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("E49").Value = "  -4.12%  "

$ws.Range("D50").Value = "0.0954"
$ws.Range("E50").Value = "  -2.30%  "

$ws.Range("D51").Value = "19.77"
$ws.Range("E51").Value = "  -7.12%  "

# Reset the style index on column D back to the default (no explicit
# style), matching the original workbook formatting while keeping the
# values stored as text.
$ws.Range("D2:D51").Style = "Normal"
